# Update "想去人数" (F column) and "最低票价" (G column) figures on the
# "展览" and "全部类型" worksheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1036
$ws1.Range("F5").Value = 2820
$ws1.Range("F7").Value = 233
$ws1.Range("F9").Value = 125
$ws1.Range("F10").Value = 78
$ws1.Range("G10").Value = 55
$ws1.Range("F11").Value = 96
$ws1.Range("F12").Value = 2651
$ws1.Range("F13").Value = 844

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1036
$ws4.Range("F6").Value = 2820
$ws4.Range("F8").Value = 233
$ws4.Range("F11").Value = 125
$ws4.Range("F12").Value = 78
$ws4.Range("G12").Value = 55
$ws4.Range("F13").Value = 96
$ws4.Range("F14").Value = 2651
$ws4.Range("F15").Value = 844
